$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1602.79
$ws.Range("I32").Value = 1538.5714
$ws.Range("J32").Value = 4749.5
$ws.Range("K32").Value = 1538.5714
$ws.Range("L32").Value = 4749.5
$ws.Range("M32").Value = -1251.5714
$ws.Range("N32").Value = -5323.5

$ws.Range("H61").Value = 11270.294
$ws.Range("I61").Value = 4633.778
$ws.Range("J61").Value = 18736.375
$ws.Range("K61").Value = 4633.778
$ws.Range("L61").Value = 18736.375
$ws.Range("M61").Value = -4421.778
$ws.Range("N61").Value = -19160.375

$ws.Range("H74").Value = 4069.7693
$ws.Range("I74").Value = 3151.8333
$ws.Range("J74").Value = 6135.125
$ws.Range("K74").Value = 3151.8333
$ws.Range("L74").Value = 6135.125
$ws.Range("M74").Value = -2277.8333
$ws.Range("N74").Value = -7883.125

$ws.Range("H77").Value = 4069.7693
$ws.Range("I77").Value = 3151.8333
$ws.Range("J77").Value = 6135.125
$ws.Range("K77").Value = 15759.1665
$ws.Range("L77").Value = 30675.625
$ws.Range("M77").Value = -11391.1665
$ws.Range("N77").Value = -39411.625

$ws.Range("H132").Value = 6346.912
$ws.Range("I132").Value = 5814.911
$ws.Range("J132").Value = 8341.916999999999
$ws.Range("K132").Value = 17444.733
$ws.Range("L132").Value = 25025.751
$ws.Range("M132").Value = -14914.733
$ws.Range("N132").Value = -30085.751

$ws.Range("H136").Value = 11270.294
$ws.Range("I136").Value = 4633.778
$ws.Range("J136").Value = 18736.375
$ws.Range("K136").Value = 13901.334
$ws.Range("L136").Value = 56209.125
$ws.Range("M136").Value = -11351.334
$ws.Range("N136").Value = -61309.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2237.7144
$ws.Range("I134").Value = 2166.6365
$ws.Range("J134").Value = 2498.3333
$ws.Range("K134").Value = 6499.9095
$ws.Range("L134").Value = 7494.999899999999
$ws.Range("M134").Value = -3964.9095
$ws.Range("N134").Value = -12564.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32865.684
$ws.Range("I31").Value = 4171.75
$ws.Range("J31").Value = 64747.832
$ws.Range("K31").Value = 4171.75
$ws.Range("L31").Value = 64747.832
$ws.Range("M31").Value = -3876.75
$ws.Range("N31").Value = -65337.832

$ws.Range("H34").Value = 32865.684
$ws.Range("I34").Value = 4171.75
$ws.Range("J34").Value = 64747.832
$ws.Range("K34").Value = 4171.75
$ws.Range("L34").Value = 64747.832
$ws.Range("M34").Value = -3969.75
$ws.Range("N34").Value = -65151.832

$ws.Range("H58").Value = 3328.5
$ws.Range("I58").Value = 1689.3529
$ws.Range("J58").Value = 8901.6
$ws.Range("K58").Value = 1689.3529
$ws.Range("L58").Value = 8901.6
$ws.Range("M58").Value = -1486.3529
$ws.Range("N58").Value = -9307.6

$ws.Range("H62").Value = 29001.666
$ws.Range("I62").Value = 28999.5
$ws.Range("J62").Value = 29006
$ws.Range("K62").Value = 28999.5
$ws.Range("L62").Value = 29006
$ws.Range("M62").Value = -28375.5
$ws.Range("N62").Value = -30254

$ws.Range("H65").Value = 29001.666
$ws.Range("I65").Value = 28999.5
$ws.Range("J65").Value = 29006
$ws.Range("K65").Value = 144997.5
$ws.Range("L65").Value = 145030
$ws.Range("M65").Value = -141877.5
$ws.Range("N65").Value = -151270

$ws.Range("H99").Value = 2610.3076
$ws.Range("I99").Value = 1889.2
$ws.Range("J99").Value = 3061
$ws.Range("K99").Value = 1889.2
$ws.Range("L99").Value = 3061
$ws.Range("M99").Value = -391.2
$ws.Range("N99").Value = -6057

$ws.Range("H126").Value = 2610.3076
$ws.Range("I126").Value = 1889.2
$ws.Range("J126").Value = 3061
$ws.Range("K126").Value = 5667.6
$ws.Range("L126").Value = 9183
$ws.Range("M126").Value = -3197.6
$ws.Range("N126").Value = -14123

$ws.Range("H132").Value = 3482.7026
$ws.Range("I132").Value = 3109.3333
$ws.Range("J132").Value = 5082.857
$ws.Range("K132").Value = 9327.999899999999
$ws.Range("L132").Value = 15248.571
$ws.Range("M132").Value = -6797.999899999999
$ws.Range("N132").Value = -20308.571

$ws.Range("H134").Value = 2647.6829
$ws.Range("I134").Value = 1636.1923
$ws.Range("J134").Value = 4400.933
$ws.Range("K134").Value = 4908.5769
$ws.Range("L134").Value = 13202.799
$ws.Range("M134").Value = -2373.5769
$ws.Range("N134").Value = -18272.799

$ws.Range("H136").Value = 3328.5
$ws.Range("I136").Value = 1689.3529
$ws.Range("J136").Value = 8901.6
$ws.Range("K136").Value = 5068.0587
$ws.Range("L136").Value = 26704.8
$ws.Range("M136").Value = -2518.0587
$ws.Range("N136").Value = -31804.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2683.182
$ws.Range("I102").Value = 2683.182
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2683.182
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1061.182
$ws.Range("N102").Value = ""

$ws.Range("H107").Value = 1291
$ws.Range("I107").Value = 499.4
$ws.Range("J107").Value = 5249
$ws.Range("K107").Value = 499.4
$ws.Range("L107").Value = 5249
$ws.Range("M107").Value = 1420.6
$ws.Range("N107").Value = -9089

$ws.Range("H113").Value = 2447.889
$ws.Range("I113").Value = 2409.25
$ws.Range("J113").Value = 2525.1667
$ws.Range("K113").Value = 2409.25
$ws.Range("L113").Value = 2525.1667
$ws.Range("M113").Value = -239.25
$ws.Range("N113").Value = -6865.1667

$ws.Range("H126").Value = 4167.5557
$ws.Range("I126").Value = 1598.8
$ws.Range("J126").Value = 7378.5
$ws.Range("K126").Value = 4796.4
$ws.Range("L126").Value = 22135.5
$ws.Range("M126").Value = -2326.4
$ws.Range("N126").Value = -27075.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1502.05
$ws.Range("I55").Value = 739.6667
$ws.Range("K55").Value = 739.6667
$ws.Range("M55").Value = -566.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 25569.928
$ws.Range("I62").Value = 24998.2
$ws.Range("J62").Value = 25887.555
$ws.Range("K62").Value = 24998.2
$ws.Range("L62").Value = 25887.555
$ws.Range("M62").Value = -24374.2
$ws.Range("N62").Value = -27135.555

$ws.Range("H65").Value = 25569.928
$ws.Range("I65").Value = 24998.2
$ws.Range("J65").Value = 25887.555
$ws.Range("K65").Value = 124991
$ws.Range("L65").Value = 129437.775
$ws.Range("M65").Value = -121871
$ws.Range("N65").Value = -135677.775

$ws.Range("H81").Value = 4392.6
$ws.Range("I81").Value = 3262.4546
$ws.Range("J81").Value = 7500.5
$ws.Range("K81").Value = 6524.9092
$ws.Range("L81").Value = 15001
$ws.Range("M81").Value = -5463.9092
$ws.Range("N81").Value = -17123

$ws.Range("H84").Value = 4392.6
$ws.Range("I84").Value = 3262.4546
$ws.Range("J84").Value = 7500.5
$ws.Range("K84").Value = 32624.546
$ws.Range("L84").Value = 75005
$ws.Range("M84").Value = -27320.546
$ws.Range("N84").Value = -85613

$ws.Range("H96").Value = 4000
$ws.Range("J96").Value = 4000
$ws.Range("L96").Value = 4000
$ws.Range("N96").Value = -6746

$ws.Range("H107").Value = 1326
$ws.Range("I107").Value = 1197.92
$ws.Range("J107").Value = 1726.25
$ws.Range("K107").Value = 3593.76
$ws.Range("L107").Value = 5178.75
$ws.Range("M107").Value = -1673.76
$ws.Range("N107").Value = -9018.75

$ws.Range("H132").Value = 3120.681
$ws.Range("I132").Value = 1858.1538
$ws.Range("J132").Value = 9275.5
$ws.Range("K132").Value = 5574.4614
$ws.Range("L132").Value = 27826.5
$ws.Range("M132").Value = -3044.4614
$ws.Range("N132").Value = -32886.5
